$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1122994652406417
$ws.Range("C2").Value = 0.7112299465240641
$ws.Range("P2").Value = 0.08021390374331551
$ws.Range("S2").Value = 0.09625668449197861
$ws.Range("B3").Value = 0.007407407407407408
$ws.Range("C3").Value = 0.007407407407407408
$ws.Range("P3").Value = 0.7925925925925926
$ws.Range("S3").Value = 0.1925925925925926
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.6136363636363636
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.03240740740740741
$ws.Range("D6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.04166666666666666
$ws.Range("J6").Value = 0.2175925925925926
$ws.Range("O6").Value = 0.01851851851851852
$ws.Range("Q6").Value = 0.1296296296296296
$ws.Range("R6").Value = 0.09259259259259259
$ws.Range("S6").Value = 0.462962962962963
$ws.Range("B7").Value = 0.07547169811320754
$ws.Range("D7").Value = 0.0440251572327044
$ws.Range("F7").Value = 0.06918238993710692
$ws.Range("J7").Value = 0.03144654088050314
$ws.Range("O7").Value = 0.006289308176100629
$ws.Range("Q7").Value = 0.1509433962264151
$ws.Range("R7").Value = 0.1006289308176101
$ws.Range("S7").Value = 0.5220125786163522
$ws.Range("B8").Value = 0.0625
$ws.Range("D8").Value = 0.02717391304347826
$ws.Range("F8").Value = 0.06793478260869565
$ws.Range("J8").Value = 0.04076086956521739
$ws.Range("O8").Value = 0.01358695652173913
$ws.Range("Q8").Value = 0.1766304347826087
$ws.Range("R8").Value = 0.1195652173913044
$ws.Range("S8").Value = 0.4918478260869565
$ws.Range("B9").Value = 0.05853658536585366
$ws.Range("D9").Value = 0.00975609756097561
$ws.Range("F9").Value = 0.03414634146341464
$ws.Range("J9").Value = 0.07804878048780488
$ws.Range("O9").Value = 0.02926829268292683
$ws.Range("Q9").Value = 0.1853658536585366
$ws.Range("R9").Value = 0.0975609756097561
$ws.Range("S9").Value = 0.5073170731707317
$ws.Range("B10").Value = 0.1101694915254237
$ws.Range("D10").Value = 0.0211864406779661
$ws.Range("E10").Value = 0.001412429378531073
$ws.Range("F10").Value = 0.06073446327683616
$ws.Range("J10").Value = 0.07062146892655367
$ws.Range("O10").Value = 0.008474576271186441
$ws.Range("Q10").Value = 0.2372881355932203
$ws.Range("R10").Value = 0.08192090395480225
$ws.Range("S10").Value = 0.4081920903954802
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.02941176470588235
$ws.Range("K11").Value = 0.1848739495798319
$ws.Range("L11").Value = 0.542016806722689
$ws.Range("S11").Value = 0.1008403361344538
$ws.Range("G12").Value = 0.7692307692307693
$ws.Range("J12").Value = 0.06923076923076923
$ws.Range("K12").Value = 0.03076923076923077
$ws.Range("L12").Value = 0.02307692307692308
$ws.Range("S12").Value = 0.1076923076923077
$ws.Range("G13").Value = 0.675
$ws.Range("J13").Value = 0.075
$ws.Range("S13").Value = 0.25
$ws.Range("F15").Value = 0.04926108374384237
$ws.Range("H15").Value = 0.1576354679802956
$ws.Range("I15").Value = 0.07389162561576355
$ws.Range("J15").Value = 0.2610837438423645
$ws.Range("K15").Value = 0.03940886699507389
$ws.Range("M15").Value = 0.01970443349753695
$ws.Range("O15").Value = 0.04926108374384237
$ws.Range("S15").Value = 0.3497536945812808
$ws.Range("F16").Value = 0.03184713375796178
$ws.Range("H16").Value = 0.1273885350318471
$ws.Range("I16").Value = 0.06369426751592357
$ws.Range("J16").Value = 0.267515923566879
$ws.Range("K16").Value = 0.1082802547770701
$ws.Range("M16").Value = 0.01910828025477707
$ws.Range("O16").Value = 0.07643312101910828
$ws.Range("S16").Value = 0.3057324840764331
$ws.Range("F17").Value = 0.01627906976744186
$ws.Range("H17").Value = 0.1511627906976744
$ws.Range("I17").Value = 0.08139534883720931
$ws.Range("J17").Value = 0.2465116279069768
$ws.Range("K17").Value = 0.09534883720930233
$ws.Range("M17").Value = 0.02325581395348837
$ws.Range("O17").Value = 0.07906976744186046
$ws.Range("S17").Value = 0.3069767441860465
$ws.Range("F18").Value = 0.02617801047120419
$ws.Range("H18").Value = 0.1884816753926702
$ws.Range("I18").Value = 0.08900523560209424
$ws.Range("J18").Value = 0.3036649214659686
$ws.Range("K18").Value = 0.06806282722513089
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.07329842931937172
$ws.Range("S18").Value = 0.2408376963350785
$ws.Range("F19").Value = 0.03817733990147783
$ws.Range("H19").Value = 0.1354679802955665
$ws.Range("I19").Value = 0.07820197044334976
$ws.Range("J19").Value = 0.188423645320197
$ws.Range("K19").Value = 0.06342364532019705
$ws.Range("M19").Value = 0.01354679802955665
$ws.Range("N19").Value = 0.0006157635467980296
$ws.Range("O19").Value = 0.05726600985221675
$ws.Range("S19").Value = 0.4248768472906404
